# Update 合肥-漫展信息.xlsx: refresh "想去人数" (interested count) figures
# generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3194
$ws1.Range("F4").Value = 150
$ws1.Range("F5").Value = 21
$ws1.Range("F6").Value = 133

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3195
$ws4.Range("F8").Value = 150
$ws4.Range("F9").Value = 21
$ws4.Range("F11").Value = 133
